$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.098109722137451
$ws.Range("B1").Value = 2.928167104721069
$ws.Range("C1").Value = 4.562941074371338
$ws.Range("D1").Value = 2.621547222137451
$ws.Range("E1").Value = 2.123500347137451
